$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 14) to the table, following the same pattern
# as the existing rows (2-13).

# Copy the date cell's style (A13 -> A14) so the new date uses the same
# date number format instead of creating a brand-new style entry.
$ws.Range("A13").Copy($ws.Range("A14"))
$ws.Range("A14").Value = 42620.886458333334

$ws.Range("B14").Value = 93
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = "Random"
